$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "42.844.17"
$ws.Range("E2").Value = "  -0.68%  "
Set-TextValue "D3" "2.367.70"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "318.16"
$ws.Range("E5").Value = "  -3.08%  "
Set-TextValue "D6" "108.70"
$ws.Range("E6").Value = "  +2.71%  "
Set-TextValue "D7" "0.635"
$ws.Range("E7").Value = "  -3.03%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue "D9" "0.623"
$ws.Range("E9").Value = "  -4.80%  "
Set-TextValue "D10" "41.91"
$ws.Range("E10").Value = "  -0.69%  "
Set-TextValue "D11" "0.0928"
$ws.Range("E11").Value = "  -1.37%  "
Set-TextValue "D12" "8.55"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("E13").Value = "  -4.36%  "
Set-TextValue "D14" "0.106"
$ws.Range("E14").Value = "  -0.28%  "
Set-TextValue "D15" "16.15"
$ws.Range("E15").Value = "  -5.79%  "
Set-TextValue "D16" "2.728.85"
$ws.Range("E16").Value = "  -1.50%  "
Set-TextValue "D17" "2.353.41"
$ws.Range("E17").Value = "  -1.81%  "
Set-TextValue "D18" "42.844.26"
$ws.Range("E18").Value = "  -1.01%  "
Set-TextValue "D19" "7.66"
$ws.Range("E19").Value = "  -1.58%  "
Set-TextValue "D20" "0.0000106"
$ws.Range("E20").Value = "  -1.78%  "
Set-TextValue "D21" "76.25"
$ws.Range("E21").Value = "  -1.21%  "
Set-TextValue "D22" "3.72"
$ws.Range("E22").Value = "  -1.32%  "
Set-TextValue "D23" "256.99"
$ws.Range("E23").Value = "  -6.21%  "
Set-TextValue "D24" "2.34"
$ws.Range("E24").Value = "  -4.08%  "
Set-TextValue "D25" "9.41"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("E26").Value = "  +0.01%  "
Set-TextValue "D27" "11.46"
$ws.Range("E27").Value = "  -4.38%  "
Set-TextValue "D28" "23.05"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D30" "172.29"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D31" "36.94"
$ws.Range("E31").Value = "  -0.88%  "
Set-TextValue "D32" "0.0896"
$ws.Range("E32").Value = "  -4.82%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D33" "2.96"
$ws.Range("E33").Value = "  -6.34%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "6.06"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  +11.57%  "
Set-TextValue "D36" "0.131"
$ws.Range("E36").Value = "  -3.16%  "
Set-TextValue "D37" "4.68"
$ws.Range("E37").Value = "  -4.52%  "
Set-TextValue "D38" "0.0364"
$ws.Range("E38").Value = "  -0.31%  "
Set-TextValue "D39" "3.93"
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("E40").Value = "  -5.90%  "
Set-TextValue "D41" "0.242"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("E42").Value = "  -5.32%  "
Set-TextValue "D43" "72.03"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("E44").Value = "  +0.00%  "
Set-TextValue "D45" "12.39"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D46" "88.31"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "113.38"
$ws.Range("E47").Value = "  -8.19%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D48" "5.59"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "9.25"
$ws.Range("E49").Value = "  -1.18%  "
Set-TextValue "D50" "77.07"
$ws.Range("E50").Value = "  +6.61%  "
Set-TextValue "D51" "1.30"
$ws.Range("E51").Value = "  -1.33%  "
